# Auto-generated edit script: update cryptos list (Price & Volume(1h) columns)
# for GitHub Actions run on Sun Nov 24 11:09:37 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.Style = 'Normal'
}

Set-TextValue $ws.Range('D2') '97.629.19'
Set-TextValue $ws.Range('E2') '  -0.93%  '
Set-TextValue $ws.Range('D3') '3.369.34'
Set-TextValue $ws.Range('E3') '  +0.39%  '
Set-TextValue $ws.Range('E4') '  +0.20%  '
Set-TextValue $ws.Range('D5') '251.02'
Set-TextValue $ws.Range('E5') '  -2.33%  '
Set-TextValue $ws.Range('D6') '663.42'
Set-TextValue $ws.Range('E6') '  -0.19%  '
Set-TextValue $ws.Range('D7') '1.42'
Set-TextValue $ws.Range('E7') '  -7.33%  '
Set-TextValue $ws.Range('D8') '0.421'
Set-TextValue $ws.Range('E8') '  -11.77%  '
Set-TextValue $ws.Range('D9') '1.00'
Set-TextValue $ws.Range('E9') '  +0.06%  '
Set-TextValue $ws.Range('D10') '1.03'
Set-TextValue $ws.Range('E10') '  -4.19%  '
Set-TextValue $ws.Range('D11') '3.366.58'
Set-TextValue $ws.Range('E11') '  +0.44%  '
Set-TextValue $ws.Range('D12') '0.211'
Set-TextValue $ws.Range('E12') '  -2.56%  '
Set-TextValue $ws.Range('D13') '40.98'
Set-TextValue $ws.Range('E13') '  -3.18%  '
Set-TextValue $ws.Range('D14') '97.368.27'
Set-TextValue $ws.Range('E14') '  -0.46%  '
Set-TextValue $ws.Range('D15') '6.15'
Set-TextValue $ws.Range('E15') '  +7.90%  '
Set-TextValue $ws.Range('D16') '0.0000257'
Set-TextValue $ws.Range('E16') '  -7.03%  '
Set-TextValue $ws.Range('D17') '3.990.18'
Set-TextValue $ws.Range('E17') '  +0.26%  '
Set-TextValue $ws.Range('D18') '8.64'
Set-TextValue $ws.Range('E18') '  +12.13%  '
Set-TextValue $ws.Range('D19') '3.374.47'
Set-TextValue $ws.Range('E19') '  +0.07%  '
Set-TextValue $ws.Range('D20') '0.573'
Set-TextValue $ws.Range('E20') '  +33.11%  '
Set-TextValue $ws.Range('D21') '17.03'
Set-TextValue $ws.Range('E21') '  +1.74%  '
Set-TextValue $ws.Range('D22') '10.97'
Set-TextValue $ws.Range('E22') '  +3.65%  '
Set-TextValue $ws.Range('D23') '503.29'
Set-TextValue $ws.Range('E23') '  -4.97%  '
Set-TextValue $ws.Range('D24') '3.38'
Set-TextValue $ws.Range('E24') '  -5.74%  '
Set-TextValue $ws.Range('D25') '0.0000201'
Set-TextValue $ws.Range('E25') '  -8.19%  '
Set-TextValue $ws.Range('D26') '6.27'
Set-TextValue $ws.Range('E26') '  +1.00%  '
Set-TextValue $ws.Range('D27') '95.27'
Set-TextValue $ws.Range('E27') '  -6.84%  '
Set-TextValue $ws.Range('D28') '12.26'
Set-TextValue $ws.Range('E28') '  -2.32%  '
Set-TextValue $ws.Range('D29') '3.559.97'
Set-TextValue $ws.Range('E29') '  +0.78%  '
Set-TextValue $ws.Range('D30') '0.149'
Set-TextValue $ws.Range('E30') '  +0.69%  '
Set-TextValue $ws.Range('D31') '11.27'
Set-TextValue $ws.Range('E31') '  +1.96%  '
Set-TextValue $ws.Range('D32') '0.995'
Set-TextValue $ws.Range('E32') '  -0.35%  '
Set-TextValue $ws.Range('E33') '  +1.28%  '
Set-TextValue $ws.Range('D34') '2.57'
Set-TextValue $ws.Range('E34') '  +21.63%  '
Set-TextValue $ws.Range('E35') '  +0.32%  '
Set-TextValue $ws.Range('D36') '0.557'
Set-TextValue $ws.Range('E36') '  +4.07%  '
Set-TextValue $ws.Range('D37') '28.76'
Set-TextValue $ws.Range('E37') '  -2.20%  '
Set-TextValue $ws.Range('D38') '7.81'
Set-TextValue $ws.Range('E38') '  -0.04%  '
Set-TextValue $ws.Range('D39') '1.47'
Set-TextValue $ws.Range('E39') '  +10.37%  '
Set-TextValue $ws.Range('D40') '526.22'
Set-TextValue $ws.Range('E40') '  +0.10%  '
Set-TextValue $ws.Range('E41') '  +0.01%  '
Set-TextValue $ws.Range('D42') '0.150'
Set-TextValue $ws.Range('E42') '  -4.44%  '
Set-TextValue $ws.Range('D43') '24.67'
Set-TextValue $ws.Range('E43') '  -0.14%  '
Set-TextValue $ws.Range('D44') '8.90'
Set-TextValue $ws.Range('E44') '  +14.48%  '
Set-TextValue $ws.Range('D45') '0.848'
Set-TextValue $ws.Range('E45') '  +3.18%  '
Set-TextValue $ws.Range('D46') '0.0424'
Set-TextValue $ws.Range('E46') '  -2.78%  '
Set-TextValue $ws.Range('D47') '3.69'
Set-TextValue $ws.Range('E47') '  -5.86%  '
Set-TextValue $ws.Range('D48') '5.67'
Set-TextValue $ws.Range('E48') '  +10.39%  '
Set-TextValue $ws.Range('D49') '1.68'
Set-TextValue $ws.Range('E49') '  +9.12%  '
Set-TextValue $ws.Range('D50') '54.10'
Set-TextValue $ws.Range('E50') '  +6.17%  '
Set-TextValue $ws.Range('D51') '3.17'
Set-TextValue $ws.Range('E51') '  -7.64%  '
